$wb = $excel.ActiveWorkbook

# --- prop_mat sheet: selection moved from C5 to B2 ---
$wsPropMat = $wb.Worksheets.Item("prop_mat")
$wsPropMat.Range("B2").Select()

# --- Add the new "varios" worksheet at the end of the workbook ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "varios"

# Match the page setup used by the rest of the workbook's sheets
$wsTemplate = $wb.Worksheets.Item(1)
$ws.PageSetup.LeftMargin = $wsTemplate.PageSetup.LeftMargin
$ws.PageSetup.RightMargin = $wsTemplate.PageSetup.RightMargin
$ws.PageSetup.TopMargin = $wsTemplate.PageSetup.TopMargin
$ws.PageSetup.BottomMargin = $wsTemplate.PageSetup.BottomMargin
$ws.PageSetup.HeaderMargin = $wsTemplate.PageSetup.HeaderMargin
$ws.PageSetup.FooterMargin = $wsTemplate.PageSetup.FooterMargin
$ws.PageSetup.CenterHeader = $wsTemplate.PageSetup.CenterHeader
$ws.PageSetup.CenterFooter = $wsTemplate.PageSetup.CenterFooter
$ws.PageSetup.PrintGridlines = $false
$ws.PageSetup.PrintHeadings = $false
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Zoom = 100

# --- Populate the cells (row by row, left to right, so shared strings are
#     appended in the same order as the reference workbook) ---

# Row 1
$ws.Range("A1").Value = "E"
$ws.Range("B1").Value = 200000000000
$ws.Range("B1").NumberFormat = "0.00E+00"
$ws.Range("C1").Value = "Pa"
$ws.Range("D1").Value = "módulo de Young"

# Row 2
$ws.Range("A2").Value = "nu"
$ws.Range("B2").Value = 0.3
$ws.Range("D2").Value = "coeficiente de Poisson"
$ws.Range("G2").Value = "Unidades de fuerza en N"

# Row 3
$ws.Range("A3").Value = "rho"
$ws.Range("B3").Formula = "=7850"
$ws.Range("C3").Value = "kg/m³"
$ws.Range("D3").Value = "densidad"
$ws.Range("G3").Value = "Unidades de longitud en m"

# Row 4
$ws.Range("A4").Value = "g"
$ws.Range("B4").Value = 9.81
$ws.Range("C4").Value = "m/s²"
$ws.Range("D4").Value = "aceleracion de la gravedad"

# Row 5
$ws.Range("A5").Value = "espesor"
$ws.Range("B5").Value = 0.01
$ws.Range("C5").Value = "m"

# Row 6
$ws.Range("A6").Value = "U_LONG"
$ws.Range("B6").Value = "m"

# Row 7
$ws.Range("A7").Value = "U_FUERZA"
$ws.Range("B7").Value = "N"

# Row 8
$ws.Range("A8").Value = "U_ESFUER"
$ws.Range("B8").Value = "Pa"

# Row 9
$ws.Range("A9").Value = "ESC_UV"
$ws.Range("B9").Value = 10000

# Row 21 - stray formatted (scientific) empty cell, matching reference file
$ws.Range("C21").NumberFormat = "0.00E+00"

$ws.Range("B5").Select()

# --- Defined names pointing at the new sheet ---
$wb.Names.Add("espesor", "=varios!`$B`$5")
$wb.Names.Add("g", "=varios!`$B`$4")
$wb.Names.Add("Poisson", "=varios!`$B`$2")
$wb.Names.Add("rho", "=varios!`$B`$3")
$wb.Names.Add("U_ESFUERZO", "=varios!`$B`$8")
$wb.Names.Add("U_FUERZA", "=varios!`$B`$7")
$wb.Names.Add("U_LONG", "=varios!`$B`$6")
$wb.Names.Add("Young", "=varios!`$B`$1")
